# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig -
# rebrand StructureDefinition spreadsheet from ibm.com/Alvearie to
# linuxforhealth.org/LinuxForHealth, bump version/date, and refresh the
# generated "Elements" constraint placement.

$wb = $excel.ActiveWorkbook

# ---- "Metadata" sheet --------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sex-assigned-at-birth"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- "Elements" sheet ---------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# The ele-1/ext-1 invariant is no longer reported against the root
# "Extension" row (row 2) ...
$elem.Range("AI2").Value = ""

# ... the fixed-value URL on Extension.url (row 5) now points at
# linuxforhealth.org ...
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sex-assigned-at-birth"

# ... and the extensible binding's ValueSet (row 7) does too.
$elem.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/sex-assigned-at-birth"

# Column Y ("Binding Value Set") widened slightly to fit the longer URL.
$elem.Columns.Item(25).ColumnWidth = 59.3
